# eval fully connected gradients
# ---------------------------------------------------------------------------
# This script reproduces the "eval fully connected gradients" commit:
#  1. Flip the sign convention of the delta (error) computation in K56:M56
#     from (target - sigmoid) to (sigmoid - target), which in turn flips the
#     sign of the dependent K59:M59 values (formulas unchanged, values
#     recalc automatically).
#  2. Add a new "delta with weights" (gradient) block starting at row 61:
#       - K61 merged across K61:M61, shared-string label "delta with weights"
#       - O61 shared-string label "Sum"
#       - rows 62..73: K/L/M = $K$59/$L$59/$M$59 * A44..A55 / B44..B55 / C44..C55
#       - rows 62..73: O = SUM(K:M) of that row
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Flip K56:M56 formulas (was K53-K47, now K47-K53, etc.)
$ws.Range("K56").Formula = "=K47-K53"
$ws.Range("L56").Formula = "=L47-L53"
$ws.Range("M56").Formula = "=M47-M53"

# K59:M59 keep their original formulas (K56*K50 etc.); their cached values
# change automatically from the K56:M56 sign flip above once recalculated.

# 2. New "delta with weights" section.
# Row 61: label cell (merged K61:M61) + "Sum" header in O61.
# Copy formatting from the existing K58:M58 label block (style s="5", the
# centred/bold "title" style used for the other section headers) so the new
# header matches the sheet's existing look-and-feel exactly.
$ws.Range("K58:M58").Copy()
$ws.Range("K61:M61").PasteSpecial(-4122)
$ws.Range("K61").Value = "δ с весами"
$ws.Range("O61").Value = "Сумма"
$ws.Range("K61:M61").Merge()

# Rows 62-73: gradients = $K$59/$L$59/$M$59 * (input column A/B/C, rows 44-55)
# plus a row-wise sum in column O.
for ($i = 0; $i -lt 12; $i++) {
    $destRow = 62 + $i
    $srcRow = 44 + $i

    $ws.Range("K$destRow").Formula = "=`$K`$59*A$srcRow"
    $ws.Range("L$destRow").Formula = "=`$L`$59*B$srcRow"
    $ws.Range("M$destRow").Formula = "=`$M`$59*C$srcRow"
    $ws.Range("O$destRow").Formula = "=SUM(K${destRow}:M${destRow})"
}

# Restore the view state (scroll position + active cell) to match the
# author's saved position.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("Q63").Select()
